$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "gtk-runtime-3.8.1-i686.exe" row). This shifts rows 3-5
# up by one, so the former rows 3/4/5 (all "payload.dll") become rows 2/3/4.
$ws.Rows(2).Delete()
